$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells (Wins / Losses / Ties) with the same header
# formatting (bold, bordered, centered) used by the rest of row 1.
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) gets the same season record: 86 wins, 76 losses,
# 0 ties.
$lastRow = 48
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 86
    $ws.Cells.Item($row, 31).Value = 76
    $ws.Cells.Item($row, 32).Value = 0
}
